# "Update fuzzy matching notes"
#
# - Un-hides all the previously filtered-out rows and clears the
#   AutoFilter's column criteria (but keeps the filter dropdowns / range).
# - Re-colours a handful of rows that had been marked "still fuzzy" (orange)
#   to the "resolved" colour (green) used elsewhere in the sheet.
# - Adds a "TODO" note in column J of row 30, next to the other notes.
# - Updates the sheet's scroll position / selection to reflect where the
#   author was last working.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear the AutoFilter criteria on column C ("Done") and reveal every row
# that had been hidden by the previous filter. In Excel's object model,
# re-applying AutoFilter on the already-filtered field toggles the
# criteria off (and shows all rows again) while leaving the filter
# dropdown buttons / range in place.
$ws.AutoFilter.Range.AutoFilter(3)

# Belt-and-braces: make sure none of the data rows are left hidden, in
# case any were hidden independently of the AutoFilter criteria.
$ws.Rows.Item("2:51").Hidden = $false

# Rows that move from "orange" (still needs review) to "green" (resolved)
# formatting - matches the colour already used on rows such as 2-12.
$resolvedRows = @(13, 15, 16, 43, 45, 46, 48, 49)
$greenColor = $ws.Range("A2").Interior.Color
foreach ($r in $resolvedRows) {
    $rng = $ws.Range("A" + $r + ":C" + $r)
    $rng.Interior.Color = $greenColor
}

# New note cell.
$ws.Range("J30").Value = "TODO"

# Reflect the author's final selection / scroll position.
$ws.Range("J31").Select()
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
